# Add two new columns, I ("I0") and J ("IF"), to the results table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers -----------------------------------------------------------
# Copy the header formatting (bold font, borders, centered alignment)
# from the existing "IP" header (H1) onto the two new header cells so
# they match the look of the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data ----------------------------------------------------------------
# Values for the new I0 / IF columns, one entry per data row (rows 2-35).
$iVals = @(2,6,6,3,1,8,6,6,3,5,9,9,5,8,5,8,7,7,7,6,7,8,6,7,8,6,7,7,7,6,3,7,6,3)
$jVals = @(5,9,6,5,1,8,6,7,4,5,9,9,5,8,6,8,8,7,7,6,8,8,6,7,8,6,7,8,7,6,4,7,7,3)

for ($k = 0; $k -lt $iVals.Count; $k++) {
    $r = $k + 2
    $ws.Cells.Item($r, 9).Value  = $iVals[$k]
    $ws.Cells.Item($r, 10).Value = $jVals[$k]
}
